$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the three existing tasks that were completed (strikethrough formatting
# reused from the existing "done" style already used elsewhere in the sheet).
$ws.Range("A16:A18").Font.Strikethrough = $true

# New tasks captured in this commit.
$ws.Range("A19").Value = "footer: data update, sources"
$ws.Range("A20").Value = "info panel : census data"
$ws.Range("A21").Value = "mess with changing symbology of waste sites: no color"
$ws.Range("A22").Value = 'description of the site after "Welcome to Haz Mat Mapper"'
$ws.Range("A23").Value = "add show all exporters/importers checkboxes to filter menu"
$ws.Range("A24").Value = "create color key to coordinate site color across icicle, map, and pov chart"

# These new rows pick up a distinct (but visually "normal") cell style - touch
# the font explicitly so a dedicated style record is written for them instead
# of silently reusing the default.
$ws.Range("A19:A24").Font.ThemeColor = 1

# A handful of blank placeholder rows were appended below the new tasks,
# carrying the same strikethrough-capable style used for completed rows.
$ws.Range("A25:A29").Font.Strikethrough = $true

# Leave the sheet scrolled/selected the way the author left it.
$ws.Range("A24").Select() | Out-Null
